$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths ---
$ws.Columns(2).ColumnWidth = 17.33   # column B -> ~18.14 chars stored
$ws.Columns(6).ColumnWidth = 13.6    # column F -> ~14.57 chars stored

# --- Update date values in column B ---
$ws.Range("B4").Value2 = 44628
$ws.Range("B5").Value2 = 44622
$ws.Range("B6").Value2 = 44623
$ws.Range("B7").Value2 = 44497
$ws.Range("B8").Value2 = 44502
$ws.Range("B9").Value2 = 44539

# --- Update CDC values in column E ---
$ws.Range("E4").Value2 = 7617
$ws.Range("E5").Value2 = 7528
$ws.Range("E6").Value2 = 7600
$ws.Range("E7").Value2 = 7610
$ws.Range("E8").Value2 = 7493
$ws.Range("E9").Value2 = 7491
$ws.Range("E10").Value2 = 7525

# --- New custom date display format for the date columns ---
$ws.Range("B3:B10").NumberFormat = "ddd dd/mm/yyyy"
$ws.Range("F3:F10").NumberFormat = "ddd dd/mm/yyyy"

# --- New labels with colored fills (order matters for shared-string index) ---
function MyRGB($r,$g,$b) { return $r + ($g*256) + ($b*65536) }

$ws.Range("A6").Value2 = "rep"
$ws.Range("A6").Interior.Color = MyRGB 0 176 80      # green FF00B050

$ws.Range("A4").Value2 = "prg"
$ws.Range("A4").Interior.Color = MyRGB 255 192 0     # orange FFFFC000

# --- Empty cell G10 keeps a quote-prefix style only ---
$ws.Range("G10").Value2 = "'"
$ws.Range("G10").ClearContents()

# --- Selection marker ---
$ws.Range("A4").Select()
